$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-06 Sunday" "2025-04-07 Monday"

Replace-Text "766×4=" "635×7="
Replace-Text "811×7=" "258×9="
Replace-Text "801×9=" "289×5="
Replace-Text "345×9=" "272×8="
Replace-Text "572×2=" "826×4="
Replace-Text "904×7=" "116×4="
Replace-Text "332×5=" "267×6="
Replace-Text "575×7=" "419×4="
Replace-Text "914×4=" "900×5="
Replace-Text "704×9=" "445×3="
Replace-Text "420×3=" "136×4="
Replace-Text "866×4=" "826×7="
Replace-Text "624×5=" "862×4="
Replace-Text "251×2=" "108×9="
Replace-Text "352×3=" "572×8="
Replace-Text "122×3=" "569×4="
Replace-Text "432×5=" "957×5="
Replace-Text "883×2=" "303×6="
Replace-Text "765×9=" "628×7="
Replace-Text "308×5=" "773×7="
Replace-Text "378×4=" "645×3="
Replace-Text "684×6=" "275×9="
Replace-Text "359×3=" "691×3="
Replace-Text "899×6=" "652×2="
Replace-Text "230×9=" "481×5="
